$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column D
$ws.Range("D1").Value = "Error"

# Update data rows 2-7 with new secant-method iteration values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = -3
$ws.Range("C2").Value = 1.53703703703704
$ws.Range("D2").Value = 1.000005

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = -2
$ws.Range("C3").Value = -0.75
$ws.Range("D3").Value = 1.000005

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = -2.32793522267206
$ws.Range("C4").Value = -0.0402948543281054
$ws.Range("D4").Value = 0.327935222672065

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = -2.3465543662801
$ws.Range("C5").Value = 0.0015198464829762
$ws.Range("D5").Value = 0.0186191436080341

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = -2.34587761290438
$ws.Range("C6").Value = -0.0000024874291515431
$ws.Range("D6").Value = 0.0006767533757217001

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = -2.34587871869074
$ws.Range("C7").Value = -0.000000000151896273337115
$ws.Range("D7").Value = 0.00000110578635981184

# Remove the now-obsolete rows 8-10 (table shrank from 10 rows to 7)
$ws.Range("A8:D10").Clear()
